$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hours")

# --- Update the window view (bookViews / workbookView) ---
$win = $excel.ActiveWindow
$win.Left = 5928
$win.Top = 8160
$win.Width = 17280
$win.Height = 9960

# --- Set column widths for Start date / End date (columns E:F) ---
$ws.Range("E:F").ColumnWidth = 10.33203125

# --- Row 2 data update (order matters for shared-string table layout) ---
$ws.Range("B2").Value = "苏鹏程"
$ws.Range("D2").Value = "Y24M11D27-0004"
$ws.Range("A2").Value = "WW-20241126002"
$ws.Range("C2").Value = "EFB"
$ws.Range("E2").Value = 45635
$ws.Range("F2").Value = 45653
$ws.Range("G2").Value = "曾祥青"
$ws.Range("H2").Value = "窦林"
$ws.Range("I2").Value = 360
$ws.Range("J2").Value = "减震器匹配，大约要更换减震器50次"

# row height
$ws.Range("A2:J2").RowHeight = 15.6

# --- Apply the new font to A2:D2 (Microsoft YaHei, black) ---
$font = $ws.Range("A2:D2").Font
$font.Name = "Microsoft YaHei"
$font.Color = 0

# --- Update selection ---
$ws.Range("K4").Select()
